$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 28.67276066666667
$ws.Range("H2").Value = 86.018282
$ws.Range("I2").Value = 0.9474462168692853
$ws.Range("J2").Value = 0.9474462168692853
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.098778
$ws.Range("N2").Value = 6.296334
$ws.Range("O2").Value = 0.3165022962792946
$ws.Range("P2").Value = 0.3165022962792947
$ws.Range("Q2").Value = 60.17775928646532
$ws.Range("R2").Value = 541.5998335781879
$ws.Range("S2").Value = 0.2998689032402593
$ws.Range("T2").Value = 0.2998689032402594

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 28.67276066666667
$ws.Range("H3").Value = 86.018282
$ws.Range("I3").Value = 0.9474462168692853
$ws.Range("J3").Value = 0.9474462168692853
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.584366666666666
$ws.Range("N3").Value = 10.7531
$ws.Range("O3").Value = 0.5405337204349202
$ws.Range("P3").Value = 0.5405337204349203
$ws.Range("Q3").Value = 102.7736875749111
$ws.Range("R3").Value = 924.9631881742
$ws.Range("S3").Value = 0.512126628516345
$ws.Range("T3").Value = 0.5121266285163452

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 28.67276066666667
$ws.Range("H4").Value = 86.018282
$ws.Range("I4").Value = 0.9474462168692853
$ws.Range("J4").Value = 0.9474462168692853
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.9480173333333334
$ws.Range("N4").Value = 2.844052
$ws.Range("O4").Value = 0.1429639832857851
$ws.Range("P4").Value = 0.1429639832857851
$ws.Range("Q4").Value = 27.18227410651822
$ws.Range("R4").Value = 244.640466958664
$ws.Range("S4").Value = 0.1354506851126808
$ws.Range("T4").Value = 0.1354506851126808

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.7569533333333333
$ws.Range("H5").Value = 2.27086
$ws.Range("I5").Value = 0.02501233070476559
$ws.Range("J5").Value = 0.02501233070476559
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.098778
$ws.Range("N5").Value = 6.296334
$ws.Range("O5").Value = 0.3165022962792946
$ws.Range("P5").Value = 0.3165022962792947
$ws.Range("Q5").Value = 1.588677003026666
$ws.Range("R5").Value = 14.29809302724
$ws.Range("S5").Value = 0.007916460103355416
$ws.Range("T5").Value = 0.007916460103355418

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.7569533333333333
$ws.Range("H6").Value = 2.27086
$ws.Range("I6").Value = 0.02501233070476559
$ws.Range("J6").Value = 0.02501233070476559
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.584366666666666
$ws.Range("N6").Value = 10.7531
$ws.Range("O6").Value = 0.5405337204349202
$ws.Range("P6").Value = 0.5405337204349203
$ws.Range("Q6").Value = 2.713198296222222
$ws.Range("R6").Value = 24.418784666
$ws.Range("S6").Value = 0.01352000817259553
$ws.Range("T6").Value = 0.01352000817259554

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.7569533333333333
$ws.Range("H7").Value = 2.27086
$ws.Range("I7").Value = 0.02501233070476559
$ws.Range("J7").Value = 0.02501233070476559
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.9480173333333334
$ws.Range("N7").Value = 2.844052
$ws.Range("O7").Value = 0.1429639832857851
$ws.Range("P7").Value = 0.1429639832857851
$ws.Range("Q7").Value = 0.7176048805244444
$ws.Range("R7").Value = 6.45844392472
$ws.Range("S7").Value = 0.003575862428814637
$ws.Range("T7").Value = 0.003575862428814637

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.8334926666666668
$ws.Range("H8").Value = 2.500478
$ws.Range("I8").Value = 0.02754145242594914
$ws.Range("J8").Value = 0.02754145242594913
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.098778
$ws.Range("N8").Value = 6.296334
$ws.Range("O8").Value = 0.3165022962792946
$ws.Range("P8").Value = 0.3165022962792947
$ws.Range("Q8").Value = 1.749316071961333
$ws.Range("R8").Value = 15.743844647652
$ws.Range("S8").Value = 0.008716932935679851
$ws.Range("T8").Value = 0.008716932935679851

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.8334926666666668
$ws.Range("H9").Value = 2.500478
$ws.Range("I9").Value = 0.02754145242594914
$ws.Range("J9").Value = 0.02754145242594913
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.584366666666666
$ws.Range("N9").Value = 10.7531
$ws.Range("O9").Value = 0.5405337204349202
$ws.Range("P9").Value = 0.5405337204349203
$ws.Range("Q9").Value = 2.987543331311111
$ws.Range("R9").Value = 26.8878899818
$ws.Range("S9").Value = 0.01488708374597965
$ws.Range("T9").Value = 0.01488708374597965

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.8334926666666668
$ws.Range("H10").Value = 2.500478
$ws.Range("I10").Value = 0.02754145242594914
$ws.Range("J10").Value = 0.02754145242594913
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.9480173333333334
$ws.Range("N10").Value = 2.844052
$ws.Range("O10").Value = 0.1429639832857851
$ws.Range("P10").Value = 0.1429639832857851
$ws.Range("Q10").Value = 0.7901654952062224
$ws.Range("R10").Value = 7.111489456856001
$ws.Range("S10").Value = 0.003937435744289638
$ws.Range("T10").Value = 0.003937435744289638
